# Ajustando o codigo para poder converter na importacao das tabelas E670LOC e E670RAT
# - Troca o conteudo da celula C4 para "TROCADO TESTE" (cria nova shared string)
# - Atualiza a selecao ativa da planilha para D10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "TROCADO TESTE"

$ws.Range("D10").Select()
